$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = "RTOS - Real Time Operating System"
$ws.Range("A7").Value = "BLDCM - Brushless DC Motor"
$ws.Range("A8").Value = "DC - Direct Current"

$ws.Range("A8").Select()
